# Apply "found mistakes were fixed" edits to sandglass_test.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New target data for column A (names) and column B (bool values), rows 2..16
$names = @(
    "Sandglass_100",
    "sandglass_81",
    "sandglass_83",
    "sandglass_84",
    "sandglass_88",
    "sandglass_89",
    "sandglass_90",
    "sandglass_91",
    "sandglass_92",
    "sandglass_93",
    "sandglass_94",
    "sandglass_95",
    "sandglass_96",
    "sandglass_97",
    "sandglass_99"
)

$values = @(1, 1, 1, 0, 0, 1, 1, 0, 0, 1, 1, 1, 1, 1, 1)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Rows 17-20 no longer hold any data - clear their contents
$ws.Range("A17:B20").ClearContents()

# Update the selection to match the saved view state
$ws.Range("A16:B16").Select()
